# Fix check for dimension level codes and text and clean up a bunch of stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell A3 should use the same "umar, eurostat" source text as A2
# (previously held the now-removed duplicate "umar eurostat" string)
$ws.Range("A3").Value = "umar, eurostat"

# Update the active selection on the sheet from A5 to A3
$ws.Range("A3").Select()

# Update the workbook window position/size
$win = $excel.ActiveWindow
$win.Left = 210
$win.Top = 150
$win.Width = 25170
$win.Height = 15300

$wb.Save()
